$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{}
$data[2] = @(2.0, 0.6666666666666666, 0.259924, 0.779772, 0.0006491957374851489, 0.0006491957374851488, 1.0, 0.3333333333333333, 0.02154, 0.06462, 0.1321623744232468, 0.1321623744232468, 0.00559876296, 0.05038886664, 0.0000857992501314881, 0.00008579925013148808)
$data[3] = @(2.0, 0.6666666666666666, 0.259924, 0.779772, 0.0006491957374851489, 0.0006491957374851488, 2.0, 0.6666666666666666, 0.08345333333333332, 0.25036, 0.5120422788703818, 0.5120422788703818, 0.02169152421333333, 0.19522371792, 0.0003324156648548338, 0.0003324156648548337)
$data[4] = @(2.0, 0.6666666666666666, 0.259924, 0.779772, 0.0006491957374851489, 0.0006491957374851488, 3.0, 1.0, 0.057988, 0.173964, 0.3557953467063713, 0.3557953467063713, 0.015072472912, 0.135652256208, 0.000230980822498827, 0.0002309808224988269)
$data[5] = @(3.0, 1.0, 383.1307676666667, 1149.392303, 0.9569214896224009, 0.9569214896224006, 1.0, 0.3333333333333333, 0.02154, 0.06462, 0.1321623744232468, 0.1321623744232468, 8.252636735540001, 74.27373061986, 0.1264690162051268, 0.1264690162051268)
$data[6] = @(3.0, 1.0, 383.1307676666667, 1149.392303, 0.9569214896224009, 0.9569214896224006, 2.0, 0.6666666666666666, 0.08345333333333332, 0.25036, 0.5120422788703818, 0.5120422788703818, 31.97353966434222, 287.76185697908, 0.4899842602462945, 0.4899842602462944)
$data[7] = @(3.0, 1.0, 383.1307676666667, 1149.392303, 0.9569214896224009, 0.9569214896224006, 3.0, 1.0, 0.057988, 0.173964, 0.3557953467063713, 0.3557953467063713, 22.21698695545467, 199.952882599092, 0.3404682131709794, 0.3404682131709793)
$data[8] = @(3.0, 1.0, 16.98778433333333, 50.963353, 0.042429314640114, 0.04242931464011399, 1.0, 0.3333333333333333, 0.02154, 0.06462, 0.1321623744232468, 0.1321623744232468, 0.36591687454, 3.29325187086, 0.005607558967988494, 0.005607558967988494)
$data[9] = @(3.0, 1.0, 16.98778433333333, 50.963353, 0.042429314640114, 0.04242931464011399, 2.0, 0.6666666666666666, 0.08345333333333332, 0.25036, 0.5120422788703818, 0.5120422788703818, 1.417687228564444, 12.75918505708, 0.02172560295923243, 0.02172560295923242)
$data[10] = @(3.0, 1.0, 16.98778433333333, 50.963353, 0.042429314640114, 0.04242931464011399, 3.0, 1.0, 0.057988, 0.173964, 0.3557953467063713, 0.3557953467063713, 0.9850876379213335, 8.865788741292, 0.01509615271289308, 0.01509615271289307)

foreach ($r in 2..10) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = 5 + $i  # column E is 5
        $ws.Cells.Item($r, $col).Value = $rowVals[$i]
    }
}

Write-Output "Done updating TPM values"
